# Generate Report for Handoff
# Updates the "Latest Handoff Date/Datetime" columns for files that were
# just handed off (status "Handback transform failed" or "Ready for handoff")
# on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$handoffDate     = "2016-03-31 08:07:29"
$handoffDatetime = "2016-03-31 08:07:17"

# Overview sheet: column D = "Latest Handoff Date"
$overview.Range("D4").Value  = $handoffDate
$overview.Range("D6").Value  = $handoffDate
$overview.Range("D7").Value  = $handoffDate
$overview.Range("D8").Value  = $handoffDate
$overview.Range("D9").Value  = $handoffDate
$overview.Range("D10").Value = $handoffDate

# zh-cn sheet: column E = "Latest Handoff Datetime"
$zhcn.Range("E4").Value  = $handoffDatetime
$zhcn.Range("E6").Value  = $handoffDatetime
$zhcn.Range("E7").Value  = $handoffDatetime
$zhcn.Range("E8").Value  = $handoffDatetime
$zhcn.Range("E9").Value  = $handoffDatetime
$zhcn.Range("E10").Value = $handoffDatetime

# de-de sheet: column E = "Latest Handoff Datetime"
$dede.Range("E4").Value  = $handoffDate
$dede.Range("E6").Value  = $handoffDate
$dede.Range("E7").Value  = $handoffDate
$dede.Range("E8").Value  = $handoffDate
$dede.Range("E9").Value  = $handoffDate
$dede.Range("E10").Value = $handoffDate
